# Commit: "Raw and Clean Data from SSA for June 14th"
# Adds the June 14 (2020-06-14, serial 43996) observation row to the daily
# sheets (out_vars, dates_dx, dates_sx, dates_deaths) and fills in the
# matching new "O" column (date 43996) in control_obs, including extending
# the shared SUM() formula in the TOTAL row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# out_vars: new row 15
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("out_vars")
$ws1.Range("A14:J14").Copy($ws1.Range("A15:J15"))
$ws1.Range("A15").Value = 43996
$ws1.Range("B15").Value = 146837
$ws1.Range("C15").Value = 207076
$ws1.Range("D15").Value = 52636
$ws1.Range("E15").Value = 17141
$ws1.Range("F15").Value = 32.503388110626069
$ws1.Range("G15").Value = 47727
$ws1.Range("H15").Value = 4323
$ws1.Range("I15").Value = 4483
$ws1.Range("J15").Value = 406549

# ---------------------------------------------------------------------
# dates_dx: row 15 already exists (blank placeholder) - fill values
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("dates_dx")
$ws2.Range("A15").Value = 43996
$ws2.Range("B15").Value = 0
$ws2.Range("C15").Value = 1
$ws2.Range("D15").Value = 1
$ws2.Range("E15").Value = 1
$ws2.Range("F15").Value = 0
$ws2.Range("G15").Value = 0
$ws2.Range("H15").Value = 0
$ws2.Range("I15").Value = 4

# ---------------------------------------------------------------------
# dates_sx: new row 15 (did not exist before)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("dates_sx")
$ws3.Range("A14:L14").Copy($ws3.Range("A15:L15"))
$ws3.Range("A15").Value = 43996
$ws3.Range("B15").Value = 0
$ws3.Range("C15").Value = 1
$ws3.Range("D15").Value = 0
$ws3.Range("E15").Value = 1
$ws3.Range("F15").Value = 1
$ws3.Range("G15").Value = 1
$ws3.Range("H15").Value = 0
$ws3.Range("I15").Value = 1
$ws3.Range("J15").Value = 1
$ws3.Range("K15").Value = 0
$ws3.Range("L15").Value = 0

# ---------------------------------------------------------------------
# dates_deaths: row 15 already exists but only A15 had a (placeholder)
# style and no value; the rest of the row was entirely empty. Copy A14's
# format onto A15 (matches the "real data" style instead of the blank
# placeholder one) before setting values.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("dates_deaths")
$ws4.Range("A14").Copy($ws4.Range("A15"))
$ws4.Range("A15").Value = 43996
$ws4.Range("B15").Value = 0
$ws4.Range("C15").Value = 0
$ws4.Range("D15").Value = 2
$ws4.Range("E15").Value = 1
$ws4.Range("F15").Value = 1
$ws4.Range("G15").Value = 1
$ws4.Range("H15").Value = 2

# ---------------------------------------------------------------------
# control_obs: fill column O (2020-06-14, serial 43996)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("control_obs")
$ws5.Range("O1").Value = 43996
$ws5.Range("O2").Value = 3465
$ws5.Range("O3").Value = 3280
$ws5.Range("O4").Value = 3280
$ws5.Range("O5").Value = 3280
$ws5.Range("O6").Value = 3280
$ws5.Range("O7").Value = 2493
$ws5.Range("O8").Value = 5102
$ws5.Range("O10").Value = 154
$ws5.Range("O11").Value = 154
$ws5.Range("O12").Value = 154
$ws5.Range("O13").Value = 154
$ws5.Range("O14").Value = 154
$ws5.Range("O15").Value = 101
$ws5.Range("O16").Value = 166
$ws5.Range("O18").Value = 807
$ws5.Range("O20").Formula = "=SUM(O2:O18)"

# ---------------------------------------------------------------------
# Selections: replicate the final cursor position recorded in each sheet
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A15").Select()

$ws2.Activate()
$ws2.Range("A15").Select()

$ws3.Activate()
$ws3.Range("A15").Select()

$ws4.Activate()
$ws4.Range("A15").Select()

$ws5.Activate()
$ws5.Range("Y15").Select()

$wb.Save()
